$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (10 and 11) -------------------------------------------
# These append new shared strings in the same left-to-right, row-by-row
# order the original author typed them in, and reuse the "ورد پرس" string
# for both D10 and D11 (matching the shared-strings table in the diff).
$ws.Range("A10").Value = "برنامه نویسی"
$ws.Range("B10").Value = "طراحی سایت"
$ws.Range("C10").Value = "چند درصد وب سایت های دنیا توسط وردپرس توسعه داده شدند؟"
$ws.Range("D10").Value = "ورد پرس"

$ws.Range("A11").Value = "دوآپس"
$ws.Range("B11").Value = "CI/CD"
$ws.Range("C11").Value = "چند درصد وب سایت های دنیا توسط وردپرس توسعه داده شدند؟ 1"
$ws.Range("D11").Value = "ورد پرس"

# --- Column width changes (B & C get wider, no longer auto bestFit) ------
$ws.Columns.Item(2).ColumnWidth = 12.8333333333
$ws.Columns.Item(3).ColumnWidth = 18.6666666667

# --- Sheet view changes: zoom in and move the selection ------------------
$excel.ActiveWindow.Zoom = 190
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C12").Select() | Out-Null
